# Update countries & provincias Spain
#
# 1) Swap the display order/content of "Fiyi" (row 205) and "Dominica" (row 206)
#    so that row 205 shows "Dominica" and row 206 shows "Fiyi".
# 2) Update a batch of per-country statistics (rows 4, 7, 18, 19, 35, 72, 104,
#    107, 170) in columns B, C, D, E, G, H.
# 3) Update the "last updated" timestamp string in cell A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Swap Fiyi / Dominica (row 205 <-> row 206, column A)
$ws.Range("A205").Value = "Dominica"
$ws.Range("A206").Value = "Fiyi"

# 2) Update statistics
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 2806513
$ws.Range("C4").Value = 26560
$ws.Range("E4").Value = 1498309
$ws.Range("G4").Value = 338
$ws.Range("H4").Value = 131136

# Row 7 - India
$ws.Range("B7").Value = 626591
$ws.Range("C7").Value = 21371
$ws.Range("E7").Value = 228579

# Row 18 - Alemania
$ws.Range("B18").Value = 196588
$ws.Range("C18").Value = 264
$ws.Range("E18").Value = 7225
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 9063

# Row 19 - Francia
$ws.Range("B19").Value = 166378
$ws.Range("C19").Value = 659
$ws.Range("D19").Value = 76802
$ws.Range("E19").Value = 59701
$ws.Range("G19").Value = 14
$ws.Range("H19").Value = 29875

# Row 35 - Emiratos Arabes Unidos
$ws.Range("B35").Value = 49469
$ws.Range("C35").Value = 400
$ws.Range("D35").Value = 38664
$ws.Range("E35").Value = 10488
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 317

# Row 72 - Uzbekistan
$ws.Range("D72").Value = 5892
$ws.Range("E72").Value = 3077

# Row 104 - Mayotte
$ws.Range("B104").Value = 2650
$ws.Range("C104").Value = 7
$ws.Range("E104").Value = 274

# Row 107 - Maldivas
$ws.Range("B107").Value = 2400
$ws.Range("C107").Value = 18
$ws.Range("D107").Value = 1969
$ws.Range("E107").Value = 421

# Row 170 - Eritrea
$ws.Range("B170").Value = 215
$ws.Range("C170").Value = 12
$ws.Range("E170").Value = 159

# 3) Update "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 2 de Julio de 2020 a las 20:43"
